$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @{ Row = 2;  Value = 4261.12301975943 },
    @{ Row = 3;  Value = 4105.214035377197 },
    @{ Row = 4;  Value = 4105.214035377197 },
    @{ Row = 5;  Value = 4105.214035377197 },
    @{ Row = 6;  Value = 4085.082942759929 },
    @{ Row = 7;  Value = 4085.082942759929 },
    @{ Row = 8;  Value = 3981.845615351559 },
    @{ Row = 9;  Value = 3981.845615351559 },
    @{ Row = 10; Value = 3981.845615351559 },
    @{ Row = 11; Value = 3981.845615351559 },
    @{ Row = 12; Value = 3981.845615351559 }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Value
}
